# Fruta / hortaliza, semanal
#
# This adds two new weekly price records for "Perejil" (Primera / Segunda)
# at Terminal Hortofrutícola Agro Chillán, inserted right before the
# existing row 52, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at row 52 (existing rows 52.. shift down to 54..)
$ws.Rows.Item(52).Insert()
$ws.Rows.Item(52).Insert()

# New row 52: Perejil, "Primera", week of 2023-06-19
$ws.Cells.Item(52, 1).Value = 7
$ws.Cells.Item(52, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52, 3).Value = "Ñuble"
$ws.Cells.Item(52, 4).Value = [datetime]"2023-06-19"
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 100112044
$ws.Cells.Item(52, 7).Value = "Perejil"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 150
$ws.Cells.Item(52, 11).Value = 1200
$ws.Cells.Item(52, 12).Value = 1200
$ws.Cells.Item(52, 13).Value = 1200
$ws.Cells.Item(52, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(52, 15).Value = "Región del Maule"
$ws.Cells.Item(52, 16).Value = 1200
$ws.Cells.Item(52, 17).Value = 1
$ws.Cells.Item(52, 18).Value = "Hortaliza"

# New row 53: Perejil, "Segunda", week of 2023-06-19
$ws.Cells.Item(53, 1).Value = 7
$ws.Cells.Item(53, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(53, 3).Value = "Ñuble"
$ws.Cells.Item(53, 4).Value = [datetime]"2023-06-19"
$ws.Cells.Item(53, 5).Value = 16
$ws.Cells.Item(53, 6).Value = 100112044
$ws.Cells.Item(53, 7).Value = "Perejil"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Segunda"
$ws.Cells.Item(53, 10).Value = 160
$ws.Cells.Item(53, 11).Value = 1000
$ws.Cells.Item(53, 12).Value = 1000
$ws.Cells.Item(53, 13).Value = 1000
$ws.Cells.Item(53, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(53, 15).Value = "Región del Maule"
$ws.Cells.Item(53, 16).Value = 1000
$ws.Cells.Item(53, 17).Value = 1
$ws.Cells.Item(53, 18).Value = "Hortaliza"
